# Refresh the cryptocurrency Price (D) and Volume(1h) (E) columns with the
# latest scraped values. Every cell in these columns is stored as literal
# text in the workbook (e.g. "67.280.77", "  -0.74%  "), so we must stop
# Excel from "helpfully" reinterpreting number-looking strings as real
# numbers when we assign them.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# These Price cells look like ordinary numbers, so force a text format on
# each of them individually before writing the new value.
$textCells = @("D5", "D6", "D11", "D13", "D15", "D20", "D21", "D24", "D25", "D28", "D29", "D35", "D40", "D42", "D44", "D45", "D46", "D47")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "67.280.77"
$ws.Range("E2").Value = "  -0.74%  "
$ws.Range("D3").Value = "3.510.01"
$ws.Range("E3").Value = "  -1.30%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "609.25"
$ws.Range("E5").Value = "  -1.33%  "
$ws.Range("D6").Value = "150.21"
$ws.Range("E6").Value = "  -2.30%  "
$ws.Range("D7").Value = "3.510.87"
$ws.Range("E7").Value = "  -1.25%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").Value = "  -1.21%  "
$ws.Range("E10").Value = "  -1.34%  "
$ws.Range("D11").Value = "7.04"
$ws.Range("E11").Value = "  +1.72%  "
$ws.Range("E12").Value = "  -1.71%  "
$ws.Range("D13").Value = "0.0000219"
$ws.Range("E13").Value = "  -2.50%  "
$ws.Range("D14").Value = "4.106.51"
$ws.Range("E14").Value = "  -1.27%  "
$ws.Range("D15").Value = "31.79"
$ws.Range("E15").Value = "  -1.07%  "
$ws.Range("D16").Value = "3.513.60"
$ws.Range("E16").Value = "  -1.28%  "
$ws.Range("D17").Value = "67.313.79"
$ws.Range("E17").Value = "  -0.76%  "
$ws.Range("E18").Value = "  +0.39%  "
$ws.Range("E19").Value = "  -0.14%  "
$ws.Range("D20").Value = "15.25"
$ws.Range("E20").Value = "  -2.71%  "
$ws.Range("D21").Value = "442.51"
$ws.Range("E21").Value = "  -2.52%  "
$ws.Range("E22").Value = "  -3.69%  "
$ws.Range("E23").Value = "  -3.25%  "
$ws.Range("D24").Value = "77.19"
$ws.Range("E24").Value = "  -0.55%  "
$ws.Range("D25").Value = "0.0000128"
$ws.Range("E25").Value = "  +9.12%  "
$ws.Range("D26").Value = "3.653.22"
$ws.Range("E26").Value = "  -1.32%  "
$ws.Range("E27").Value = "  -0.10%  "
$ws.Range("D28").Value = "10.23"
$ws.Range("E28").Value = "  -3.92%  "
$ws.Range("D29").Value = "8.34"
$ws.Range("E29").Value = "  -0.68%  "
$ws.Range("E30").Value = "  -2.32%  "
$ws.Range("E31").Value = "  -4.82%  "
$ws.Range("E32").Value = "  -0.07%  "
$ws.Range("E33").Value = "  +3.42%  "
$ws.Range("E34").Value = "  -0.78%  "
$ws.Range("D35").Value = "6.13"
$ws.Range("E35").Value = "  -1.45%  "
$ws.Range("D36").Value = "3.503.96"
$ws.Range("E36").Value = "  -1.53%  "
$ws.Range("E37").Value = "  -3.91%  "
$ws.Range("E38").Value = "  -1.23%  "
$ws.Range("E39").Value = "  +0.01%  "
$ws.Range("D40").Value = "178.23"
$ws.Range("E40").Value = "  +0.76%  "
$ws.Range("E41").Value = "  -0.04%  "
$ws.Range("D42").Value = "2.15"
$ws.Range("E42").Value = "  +2.53%  "
$ws.Range("E43").Value = "  -1.75%  "
$ws.Range("D44").Value = "5.42"
$ws.Range("E44").Value = "  -3.16%  "
$ws.Range("D45").Value = "0.879"
$ws.Range("E45").Value = "  -1.57%  "
$ws.Range("D46").Value = "45.55"
$ws.Range("E46").Value = "  -1.50%  "
$ws.Range("D47").Value = "27.69"
$ws.Range("E47").Value = "  -4.81%  "
$ws.Range("E48").Value = "  +4.47%  "
$ws.Range("E49").Value = "  -0.95%  "
$ws.Range("E51").Value = "  -1.70%  "

# Restore the default styling so no stray number-format is left on the cells.
foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}

